$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Carry over the existing header / data-row formatting (bold
#    centered header with a heavy border, plain centered data rows
#    with the same heavy border) onto the new table footprint before
#    the old content is overwritten, so the style table is reused
#    instead of growing needlessly.
# ------------------------------------------------------------------
$ws.Range("B2").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

$ws.Range("B3").Copy()
$ws.Range("A2:B6").PasteSpecial(-4122)

# The header/data border goes from "thick" to "medium".
$ws.Range("A1:B6").Borders.Weight = -4138

# ------------------------------------------------------------------
# 2. Drop the old Item/Price/Sales content that isn't part of the new
#    layout.
# ------------------------------------------------------------------
$ws.Range("C2:D5").Clear()

# ------------------------------------------------------------------
# 3. New header + student data.
# ------------------------------------------------------------------
$ws.Range("A1").Value = "Student Name "
$ws.Range("B1").Value = "Marks"

$ws.Range("A2").Value = "Riya "
$ws.Range("A3").Value = "Siya "
$ws.Range("A4").Value = "Raju "
$ws.Range("A5").Value = "Rahul"
$ws.Range("A6").Value = "Ram"

$ws.Range("B2").Value = 45
$ws.Range("B3").Value = 72
$ws.Range("B4").Value = 55
$ws.Range("B5").Value = 80
$ws.Range("B6").Value = 30

# ------------------------------------------------------------------
# 4. Pass/fail formula column - give it the plain font used elsewhere
#    in the table but with no border and default (general) alignment.
# ------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("C2:C6").PasteSpecial(-4122)
$ws.Range("C2:C6").Borders.LineStyle = -4142
$ws.Range("C2:C6").HorizontalAlignment = 1

$ws.Range("C2:C6").Formula = '=IF(B2>=40,"pass","fail")'

# ------------------------------------------------------------------
# 5. Lone formatted-but-empty cell G9 (plain font, thin border).
# ------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("G9").PasteSpecial(-4122)
$ws.Range("G9").Borders.Weight = 2
$ws.Range("G9").HorizontalAlignment = 1

# ------------------------------------------------------------------
# 6. Column widths for the new table.
# ------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 13.88
$ws.Columns("B").ColumnWidth = 6.25

Write-Output "done"
